# Update countries & provincias Spain
#
# This reproduces a refreshed data pull of the COVID dashboard ("Pais" sheet):
#   - the "last updated" timestamp moves from 10:08 to 11:25
#   - a handful of countries get newer totals (cases/new cases/active/recovered/
#     deaths-today/deaths)
#   - four country pairs swap ranking order (because one of the pair grew past
#     the other): Nepal/Austria, Nicaragua/Hong Kong, Islas Turcas y Caicos/
#     Brunei and Timor Oriental/Santa Lucia. The row that moves up gets the
#     refreshed numbers while the row that moves down keeps its previous values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Datos actualizados a ..." timestamp
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 7 de Agosto de 2020 a las 11:25"

# Estados Unidos
$ws.Cells.Item(4, 2).Value = 5032561
$ws.Cells.Item(4, 3).Value = 382
$ws.Cells.Item(4, 5).Value = 2292407

# Banglades
$ws.Cells.Item(18, 2).Value = 252502
$ws.Cells.Item(18, 3).Value = 2851
$ws.Cells.Item(18, 4).Value = 145584
$ws.Cells.Item(18, 5).Value = 103585
$ws.Cells.Item(18, 7).Value = 27
$ws.Cells.Item(18, 8).Value = 3333

# Filipinas
$ws.Cells.Item(25, 2).Value = 122754
$ws.Cells.Item(25, 3).Value = 3379
$ws.Cells.Item(25, 4).Value = 66852
$ws.Cells.Item(25, 5).Value = 53734
$ws.Cells.Item(25, 7).Value = 24
$ws.Cells.Item(25, 8).Value = 2168

# Indonesia
$ws.Cells.Item(26, 2).Value = 121226
$ws.Cells.Item(26, 3).Value = 2473
$ws.Cells.Item(26, 4).Value = 77557
$ws.Cells.Item(26, 5).Value = 38076
$ws.Cells.Item(26, 7).Value = 72
$ws.Cells.Item(26, 8).Value = 5593

# Oman
$ws.Cells.Item(35, 5).Value = 9301
$ws.Cells.Item(35, 7).Value = 10
$ws.Cells.Item(35, 8).Value = 502

# Republica de Africa Central (row 49)
$ws.Cells.Item(49, 2).Value = 50324
$ws.Cells.Item(49, 3).Value = 809
$ws.Cells.Item(49, 4).Value = 36041
$ws.Cells.Item(49, 5).Value = 12496
$ws.Cells.Item(49, 7).Value = 13
$ws.Cells.Item(49, 8).Value = 1787

# Nepal / Austria swap: Austria moves up to row 68 with refreshed numbers,
# Nepal moves down to row 69 keeping its previous numbers.
$ws.Cells.Item(68, 1).Value = "Austria"
$ws.Cells.Item(68, 2).Value = 21837
$ws.Cells.Item(68, 3).Value = 141
$ws.Cells.Item(68, 4).Value = 19690
$ws.Cells.Item(68, 5).Value = 1427
$ws.Cells.Item(68, 7).Value = 1
$ws.Cells.Item(68, 8).Value = 720

$ws.Cells.Item(69, 1).Value = "Nepal"
$ws.Cells.Item(69, 2).Value = 21750
$ws.Cells.Item(69, 4).Value = 15389
$ws.Cells.Item(69, 5).Value = 6296
$ws.Cells.Item(69, 8).Value = 65

# Row 72 (unaffected reordering, just refreshed numbers)
$ws.Cells.Item(72, 2).Value = 20272
$ws.Cells.Item(72, 3).Value = 410
$ws.Cells.Item(72, 5).Value = 8859

# Row 94
$ws.Cells.Item(94, 2).Value = 7554
$ws.Cells.Item(94, 3).Value = 22
$ws.Cells.Item(94, 5).Value = 243

# Row 102
$ws.Cells.Item(102, 2).Value = 5334
$ws.Cells.Item(102, 3).Value = 4
$ws.Cells.Item(102, 4).Value = 5066
$ws.Cells.Item(102, 5).Value = 209

# Row 110
$ws.Cells.Item(110, 2).Value = 4395
$ws.Cells.Item(110, 3).Value = 56
$ws.Cells.Item(110, 5).Value = 3034
$ws.Cells.Item(110, 7).Value = 13
$ws.Cells.Item(110, 8).Value = 97

# Nicaragua / Hong Kong swap: Hong Kong moves up to row 111 with refreshed
# numbers, Nicaragua moves down to row 112 keeping its previous numbers.
$ws.Cells.Item(111, 1).Value = "Hong Kong"
$ws.Cells.Item(111, 2).Value = 3939
$ws.Cells.Item(111, 3).Value = 89
$ws.Cells.Item(111, 4).Value = 2620
$ws.Cells.Item(111, 5).Value = 1273
$ws.Cells.Item(111, 8).Value = 46

$ws.Cells.Item(112, 1).Value = "Nicaragua"
$ws.Cells.Item(112, 2).Value = 3902
$ws.Cells.Item(112, 4).Value = 2913
$ws.Cells.Item(112, 5).Value = 866
$ws.Cells.Item(112, 8).Value = 123

# Row 124
$ws.Cells.Item(124, 2).Value = 2523
$ws.Cells.Item(124, 3).Value = 43
$ws.Cells.Item(124, 4).Value = 1846
$ws.Cells.Item(124, 5).Value = 646
$ws.Cells.Item(124, 7).Value = 2
$ws.Cells.Item(124, 8).Value = 31

# Row 126
$ws.Cells.Item(126, 2).Value = 2233
$ws.Cells.Item(126, 3).Value = 10
$ws.Cells.Item(126, 4).Value = 1927
$ws.Cells.Item(126, 5).Value = 181

# Row 127
$ws.Cells.Item(127, 2).Value = 2194
$ws.Cells.Item(127, 3).Value = 23
$ws.Cells.Item(127, 4).Value = 1658
$ws.Cells.Item(127, 5).Value = 455

# Row 134
$ws.Cells.Item(134, 2).Value = 1932
$ws.Cells.Item(134, 3).Value = 2
$ws.Cells.Item(134, 5).Value = 97

# Islas Turcas y Caicos / Brunei swap: Brunei moves up to row 185 with
# refreshed numbers, Islas Turcas y Caicos moves down to row 186 keeping its
# previous numbers.
$ws.Cells.Item(185, 1).Value = "Brunei"
$ws.Cells.Item(185, 2).Value = 142
$ws.Cells.Item(185, 3).Value = 1
$ws.Cells.Item(185, 4).Value = 138
$ws.Cells.Item(185, 5).Value = 1
$ws.Cells.Item(185, 8).Value = 3

$ws.Cells.Item(186, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(186, 3).Value = 12
$ws.Cells.Item(186, 4).Value = 39
$ws.Cells.Item(186, 5).Value = 100
$ws.Cells.Item(186, 8).Value = 2

# Timor Oriental / Santa Lucia swap: both rows already shared identical
# totals, so only the country names trade places (row 202 <-> row 203).
$ws.Cells.Item(202, 1).Value = "Santa Lucia"
$ws.Cells.Item(203, 1).Value = "Timor Oriental"
